$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 'P5: 1. SHORT VOWEL A - Pack 1'
$ws.Range("B6").Value = 'CVC words with ''a'' as in cat (Pack 1 of 2)'
$ws.Range("C6").Value = 'assistant, attack, attic, bad, bag, bat, cab, cap, cat, dab, fan, fat, gal, gap, gas, ham, has, hat, jab, jam, kit, lap, mad, man, map, mat, max, nag, nap, napkin'

$ws.Range("A7").Value = 'P6: 1. SHORT VOWEL A - Pack 2'
$ws.Range("B7").Value = 'CVC words with ''a'' as in cat (Pack 2 of 2)'
$ws.Range("C7").Value = 'pal, pan, pasta, pat, rag, ram, ran, rap, rat, sad, sap, sat, span, spank, spick, stack, tab, tact, tag, tan, tap, tax, van, wag, wax, yam, zap'

$ws.Range("A8").Value = 'P7: 1. SHORT VOWEL E - Pack 1'
$ws.Range("B8").Value = 'CVC words with ''e'' as in bed (Pack 1 of 2)'
$ws.Range("C8").Value = 'antenna, antiseptic, arrest, bed, beg, bell, ben, best, bet, cassette, deck, den, fed, fell, gem, hem, hen, insect, inspect, jet, kitten, led, leg, less, let, men, mess, met, neck, nest'

$ws.Range("A9").Value = 'P8: 1. SHORT VOWEL E - Pack 2'
$ws.Range("B9").Value = 'CVC words with ''e'' as in bed (Pack 2 of 2)'
$ws.Range("C9").Value = 'net, packet, peck, peg, pen, pest, pet, red, rest, sell, sense, sent, set, speck, step, tell, ten, tennis, tense, tent, test, ticket, vest, vet, wed, well, west, wet, yell, yet'

$ws.Range("A10").Value = 'P9: 1. SHORT VOWEL I - Pack 1'
$ws.Range("B10").Value = 'CVC words with ''i'' as in sit (Pack 1 of 2)'
$ws.Range("C10").Value = 'ant, assist, bib, bin, bit, did, dig, dim, din, dip, fig, fin, fit, fix, gig, hid, hip, hit, insist, its, jig, kick, kid, lick, lid, lip, lit, mix, nip, nit'

$ws.Range("A11").Value = 'P10: 1. SHORT VOWEL I - Pack 2'
$ws.Range("B11").Value = 'CVC words with ''i'' as in sit (Pack 2 of 2)'
$ws.Range("C11").Value = 'pant, pants, pick, pig, pin, pip, pit, rib, rid, rig, rim, rip, sick, sip, sit, six, snap, snip, spat, spin, spit, tick, tin, tip, wick, wig, win, wit, zip'

$ws.Range("A12").Value = 'P11: 1. SHORT VOWEL O - Pack 1'
$ws.Range("B12").Value = 'CVC words with ''o'' as in dog (Pack 1 of 3)'
$ws.Range("C12").Value = 'across, adopt, block, bob, bog, box, cannon, cannot, carrot, clock, cob, cod, cog, comic, connect, correct, cost, cot, cotton, crop, cross, desktop, dock, dog, dot, dragon, drop, fog, fox, god'

$ws.Range("A13").Value = 'P12: 1. SHORT VOWEL O - Pack 2'
$ws.Range("B13").Value = 'CVC words with ''o'' as in dog (Pack 2 of 3)'
$ws.Range("C13").Value = 'haddock, hog, hop, hot, incorrect, job, jog, jot, kiosk, lock, log, lot, maggot, mob, mock, mom, mop, moss, nod, nonsense, odd, opinion, parrot, pocket, pod, pond, pop, pot, pox, protect'

$ws.Range("A14").Value = 'P13: 1. SHORT VOWEL O - Pack 3'
$ws.Range("B14").Value = 'CVC words with ''o'' as in dog (Pack 3 of 3)'
$ws.Range("C14").Value = 'reckon, rob, rock, rocket, rod, rot, second, shock, sob, sock, spot, spotted, stop, ticktock, top, topic, trod, trot'

$ws.Range("A15").Value = 'P14: 1. SHORT VOWEL U - Pack 1'
$ws.Range("B15").Value = 'CVC words with ''u'' as in cup (Pack 1 of 3)'
$ws.Range("C15").Value = 'bud, bug, bump, bun, bus, crust, cub, cud, cup, cut, discuss, drug, drum, drunk, duck, dug, dumb, dump, dust, eggcup, fun, grunt, gum, gun, gust, gut, hiccup, hippopotamus, hug, hum'

$ws.Range("A16").Value = 'P15: 1. SHORT VOWEL U - Pack 2'
$ws.Range("B16").Value = 'CVC words with ''u'' as in cup (Pack 2 of 3)'
$ws.Range("C16").Value = 'hump, hunt, hut, instruct, jug, jump, jut, lump, minimum, mud, mug, mumps, must, nun, nut, product, pub, pump, pumpkin, pun, pup, puppet, putt, rub, rucksack, rug, rump, run, rust, rut'

$ws.Range("A17").Value = 'P16: 1. SHORT VOWEL U - Pack 3'
$ws.Range("B17").Value = 'CVC words with ''u'' as in cup (Pack 3 of 3)'
$ws.Range("C17").Value = 'scrub, skunk, snug, spun, stuck, stump, sub, suck, sudden, sum, sun, sunk, sunset, suntan, tantrum, truck, trumpet, trunk, trust, tub, tuck, tug, tusk, undid, undress, unpack, upon, upset, us, yum'

$ws.Range("A18").Value = 'P17: 2. 3-LETTER BLENDS'
$ws.Range("B18").Value = 'scr, spr, str, spl, thr'
$ws.Range("C18").Value = 'scrape, screen, screw, splash, splendid, splint, sprain, spray, spread, spring, sprout, strap, stream, street, string, strip, strong, three, throat, throne, throw'

$ws.Range("A19").Value = 'P18: 2. L-BLENDS - Pack 1'
$ws.Range("B19").Value = 'bl, cl, fl, gl, pl, sl (Pack 1 of 2)'
$ws.Range("C19").Value = 'black, blast, blend, bless, blob, blot, blue, blunt, clamp, clap, clip, club, clump, flag, flan, flat, flip, flock, flop, flow, glad, glass, glide, glow, glue, glum, plan, plank, play, plot'

$ws.Range("A20").Value = 'P19: 2. L-BLENDS - Pack 2'
$ws.Range("B20").Value = 'bl, cl, fl, gl, pl, sl (Pack 2 of 2)'
$ws.Range("C20").Value = 'plug, plum, plump, plus, skill, skull, slam, slap, sled, slept, slid, slim, slip, slit, slot, slow, slug, smell, spell, spelt, spill, split'

$ws.Range("A21").Value = 'P20: 2. R-BLENDS - Pack 1'
$ws.Range("B21").Value = 'br, cr, dr, fr, gr, pr, tr (Pack 1 of 2)'
$ws.Range("C21").Value = 'brave, brick, bring, brown, brush, crab, crack, crash, crown, cry, drab, drag, dress, drill, drip, frank, frantic, free, fresh, frill, frog, frost, grab, gram, green, grim, grin, grip, pray, press'

Write-Output "done"